$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data: ldr -> ldr1, add ldr2 and limit switch rows ---
$ws.Range("B10").Value = "ldr1"
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "ldr2"
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "limit switch"

# --- Arm base movement test table (columns K:L) ---
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = "trail"
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 4
$ws.Range("K9").Value = 5
$ws.Range("L9").Value = "convire"
$ws.Range("K10").Value = 6
$ws.Range("K11").Value = 7

# --- Servo PIN table (columns N:T) ---
$ws.Range("O1").Value = "s1"
$ws.Range("P1").Value = "s2"
$ws.Range("Q1").Value = "s3"
$ws.Range("R1").Value = "s4"
$ws.Range("S1").Value = "s5"
$ws.Range("T1").Value = "s6"
$ws.Range("N2").Value = "PIN"
$ws.Range("N1").Value = "Servo"

$ws.Range("O2").Value = 11
$ws.Range("P2").Value = 12
$ws.Range("Q2").Value = 13
$ws.Range("R2").Value = 50
$ws.Range("S2").Value = 51
$ws.Range("T2").Value = 24

# Copy the centered data style (A2:B2) down onto the new / trailing rows
$ws.Range("A2:B2").Copy()
$ws.Range("A11:B18").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("K5:L7").Borders.LineStyle = 1
$ws.Range("K5:L7").Borders.Weight = 2
$ws.Range("K9:L11").Borders.LineStyle = 1
$ws.Range("K9:L11").Borders.Weight = 2

$ntRange = $ws.Range("N1:T2")
$ntRange.Borders.LineStyle = 1
$ntRange.Borders.Weight = 2
$ws.Range("O1:T2").HorizontalAlignment = -4108

$ws.Range("M11").Select()
